$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-28 22:18:31'
$ws.Range("O2").Value = '2.7 °C'
$ws.Range("E3").Value = '2026-02-28 22:18:34'
$ws.Range("N3").Value = '-3.3 °C 21:30 TU'
$ws.Range("O3").Value = '-1.4 °C'
$ws.Range("E4").Value = '2026-02-28 22:18:36'
$ws.Range("H4").Value = '''84%'
$ws.Range("J4").Value = '1024.9 hPa'
$ws.Range("E5").Value = '2026-02-28 22:18:39'
$ws.Range("N5").Value = '-3.4 °C 21:58 TU'
$ws.Range("E6").Value = '2026-02-28 22:18:41'
$ws.Range("J6").Value = '1024.8 hPa'
$ws.Range("N6").Value = '8.4 °C 21:53 TU'
$ws.Range("E7").Value = '2026-02-28 22:18:44'
$ws.Range("O7").Value = '13.3 °C'
$ws.Range("E8").Value = '2026-02-28 22:18:47'
$ws.Range("E9").Value = '2026-02-28 22:18:49'
$ws.Range("E10").Value = '2026-02-28 22:18:52'
$ws.Range("O10").Value = '10.6 °C'
$ws.Range("E11").Value = '2026-02-28 22:18:55'
$ws.Range("E12").Value = '2026-02-28 22:18:58'
$ws.Range("H12").Value = '''83%'
$ws.Range("E13").Value = '2026-02-28 22:19:00'
$ws.Range("J13").Value = '1024.3 hPa'
$ws.Range("E14").Value = '2026-02-28 22:19:03'
$ws.Range("H14").Value = '''84%'
$ws.Range("E15").Value = '2026-02-28 22:19:05'
$ws.Range("O15").Value = '10.8 °C'
$ws.Range("E16").Value = '2026-02-28 22:19:08'
$ws.Range("H16").Value = '''69%'
$ws.Range("N16").Value = '-3.1 °C 21:42 TU'
$ws.Range("O16").Value = '-1.3 °C'
$ws.Range("E17").Value = '2026-02-28 22:19:10'
$ws.Range("E18").Value = '2026-02-28 22:19:13'
$ws.Range("H18").Value = '''84%'
$ws.Range("N18").Value = '6.3 °C 21:58 TU'
$ws.Range("O18").Value = '11.3 °C'
$ws.Range("E19").Value = '2026-02-28 22:19:16'
$ws.Range("I19").Value = '0.6 mm'
$ws.Range("E20").Value = '2026-02-28 22:19:18'
$ws.Range("H20").Value = '''66%'
$ws.Range("N20").Value = '-2.3 °C 21:58 TU'
$ws.Range("E21").Value = '2026-02-28 22:19:21'
$ws.Range("E22").Value = '2026-02-28 22:19:23'
$ws.Range("H22").Value = '''71%'
$ws.Range("N22").Value = '-3.2 °C 21:59 TU'
$ws.Range("E23").Value = '2026-02-28 22:19:26'
$ws.Range("H23").Value = '''73%'
$ws.Range("E24").Value = '2026-02-28 22:19:29'
$ws.Range("J24").Value = '1025.3 hPa'
$ws.Range("E25").Value = '2026-02-28 22:19:31'
$ws.Range("H25").Value = '''65%'
$ws.Range("I25").Value = '1.7 mm'
$ws.Range("O25").Value = '1.0 °C'
$ws.Range("E26").Value = '2026-02-28 22:19:34'
$ws.Range("H26").Value = '''82%'
$ws.Range("I26").Value = '0.6 mm'
$ws.Range("J26").Value = '1024.5 hPa'
$ws.Range("E27").Value = '2026-02-28 22:19:37'
$ws.Range("H27").Value = '''59%'
$ws.Range("N27").Value = '-0.8 °C 21:53 TU'
$ws.Range("E28").Value = '2026-02-28 22:19:39'
$ws.Range("I28").Value = '0.1 mm'
$ws.Range("E29").Value = '2026-02-28 22:19:42'
$ws.Range("E30").Value = '2026-02-28 22:19:45'
$ws.Range("E31").Value = '2026-02-28 22:19:47'
$ws.Range("H31").Value = '''79%'
$ws.Range("J31").Value = '1024.4 hPa'
$ws.Range("E32").Value = '2026-02-28 22:19:50'
$ws.Range("E33").Value = '2026-02-28 22:19:53'
$ws.Range("E34").Value = '2026-02-28 22:19:55'
$ws.Range("H34").Value = '''71%'
$ws.Range("I34").Value = '1.6 mm'
$ws.Range("E35").Value = '2026-02-28 22:19:58'
$ws.Range("J35").Value = '1025.1 hPa'
$ws.Range("O35").Value = '5.8 °C'
$ws.Range("E36").Value = '2026-02-28 22:20:01'
$ws.Range("J36").Value = '1025.0 hPa'
$ws.Range("E37").Value = '2026-02-28 22:20:03'
$ws.Range("H37").Value = '''81%'
$ws.Range("I37").Value = '0.6 mm'
$ws.Range("E38").Value = '2026-02-28 22:20:06'
$ws.Range("N38").Value = '8.4 °C 21:51 TU'
$ws.Range("E39").Value = '2026-02-28 22:20:09'
$ws.Range("E40").Value = '2026-02-28 22:20:11'
$ws.Range("O40").Value = '7.8 °C'
$ws.Range("E41").Value = '2026-02-28 22:20:14'
$ws.Range("H41").Value = '''74%'
$ws.Range("J41").Value = '1024.7 hPa'
$ws.Range("E42").Value = '2026-02-28 22:20:17'
$ws.Range("O42").Value = '10.9 °C'
$ws.Range("E43").Value = '2026-02-28 22:20:19'
$ws.Range("O43").Value = '7.9 °C'
$ws.Range("E44").Value = '2026-02-28 22:20:22'
$ws.Range("E45").Value = '2026-02-28 22:20:24'
$ws.Range("N45").Value = '4.1 °C 21:38 TU'
$ws.Range("E46").Value = '2026-02-28 22:20:27'
